# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) in the pitching log was originally populated with the
# wrong stat (Strike# instead of strikeouts, K). This re-derives/writes the
# corrected K values for each start (row) back into column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> corrected K value, as recomputed from the source data.
$kValues = [ordered]@{
    2  = 0
    3  = 1
    4  = 3
    5  = 1
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 2
    13 = 1
    14 = 1
    15 = 4
    16 = 3
    17 = 2
    18 = 0
    19 = 2
    20 = 0
    21 = 1
    22 = 2
    23 = 1
    24 = 1
    25 = 2
    26 = 0
    27 = 1
    28 = 1
    29 = 0
    30 = 0
    31 = 2
    32 = 1
    33 = 2
    34 = 1
    35 = 0
    36 = 0
    37 = 2
    38 = 1
    39 = 1
    40 = 0
    41 = 0
    42 = 1
    43 = 2
    44 = 1
    45 = 1
    46 = 1
    47 = 1
    48 = 0
    49 = 2
    50 = 0
    51 = 1
    53 = 1
    54 = 0
    55 = 0
    57 = 0
    58 = 1
    59 = 0
    60 = 3
    61 = 2
    62 = 2
    63 = 0
    64 = 0
    65 = 0
    66 = 2
    67 = 1
    68 = 0
    70 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
